$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "University of Puerto Rico" research-assistant entries: consolidate the
# "what / department" strings into a single "Department of Chemistry" /
# " University of Puerto Rico, San Juan, PR" pairing (matches L. Abad CV).
$ws.Range("G6").Value = "Department of Chemistry"
$ws.Range("H6").Value = " University of Puerto Rico, San Juan, PR"
$ws.Range("G10").Value = "Department of Chemistry"
$ws.Range("H10").Value = " University of Puerto Rico, San Juan, PR"

# Widen column H so the long institution/department text is readable.
$ws.Columns.Item(8).ColumnWidth = 45.86

# Move the cursor/selection like the author left it.
$ws.Range("H18").Select() | Out-Null
